$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 15-27), columns B:E
$data = @(
    @(10862, 11000, 8424, 4291),
    @(3708, 9793, 4753, 6798),
    @(6041, 6504, 10994, 9278),
    @(6816, 5936, 3781, 8295),
    @(4720, 3147, 8214, 7785),
    @(4723, 5815, 4290, 4781),
    @(8683, 7161, 8814, 9272),
    @(6901, 4751, 4546, 5890),
    @(5681, 9793, 6209, 6816),
    @(5936, 10769, 4720, 6041),
    @(4290, 10862, 4721, 10994),
    @(6798, 7785, 10997, 8814),
    @(8295, 4546, 8214, 7161)
)

$startRow = 15
$endRow = 27

# Column A: extend the running-count series (=A(row-1)+1) for all new rows
# in one range assignment so the engine groups it as a shared formula,
# matching how dragging the fill handle down extends the existing series.
$ws.Range("A$startRow`:A$endRow").Formula = "=A$($startRow - 1)+1"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}

# Update selection to match the final state
$ws.Range("F27").Select() | Out-Null
